$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.99 -> 0M
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

# Row 2: 0 -> 0M
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"

# Row 3: 27 -> 0M
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (i.e. before the current row 4),
# each containing a single value, in order.
$target = $t.Rows.Item(4)
$newValues = @("45","0.00003","0.00005","0.00003","0.00000","0.00003","0.00004","0.00004","0.00146","100.0")
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($target)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# The three trailing multi-value (tab-separated) rows collapse down to a
# single value each. After the 10-row insertion above, the table now has
# 46 rows total, and these three rows are the last three (44, 45, 46).
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = "99.99"
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = "0"
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = "27"

Write-Output "done"
